$d = $word.ActiveDocument

$d.Content.Find.Execute("939×7=", $true, $false, $false, $false, $false, $true, 1, $false, "394×9=", 2) | Out-Null
$d.Content.Find.Execute("679×3=", $true, $false, $false, $false, $false, $true, 1, $false, "791×9=", 2) | Out-Null
$d.Content.Find.Execute("484×9=", $true, $false, $false, $false, $false, $true, 1, $false, "195×9=", 2) | Out-Null
$d.Content.Find.Execute("539×8=", $true, $false, $false, $false, $false, $true, 1, $false, "693×6=", 2) | Out-Null
$d.Content.Find.Execute("437×7=", $true, $false, $false, $false, $false, $true, 1, $false, "432×4=", 2) | Out-Null
$d.Content.Find.Execute("690×2=", $true, $false, $false, $false, $false, $true, 1, $false, "565×7=", 2) | Out-Null
$d.Content.Find.Execute("593×6=", $true, $false, $false, $false, $false, $true, 1, $false, "353×7=", 2) | Out-Null
$d.Content.Find.Execute("556×8=", $true, $false, $false, $false, $false, $true, 1, $false, "862×5=", 2) | Out-Null
$d.Content.Find.Execute("651×7=", $true, $false, $false, $false, $false, $true, 1, $false, "813×6=", 2) | Out-Null
$d.Content.Find.Execute("215×6=", $true, $false, $false, $false, $false, $true, 1, $false, "217×9=", 2) | Out-Null
$d.Content.Find.Execute("170×3=", $true, $false, $false, $false, $false, $true, 1, $false, "722×7=", 2) | Out-Null
$d.Content.Find.Execute("199×6=", $true, $false, $false, $false, $false, $true, 1, $false, "962×3=", 2) | Out-Null
$d.Content.Find.Execute("294×3=", $true, $false, $false, $false, $false, $true, 1, $false, "320×9=", 2) | Out-Null
$d.Content.Find.Execute("506×9=", $true, $false, $false, $false, $false, $true, 1, $false, "490×5=", 2) | Out-Null
$d.Content.Find.Execute("490×4=", $true, $false, $false, $false, $false, $true, 1, $false, "235×6=", 2) | Out-Null
$d.Content.Find.Execute("231×3=", $true, $false, $false, $false, $false, $true, 1, $false, "207×2=", 2) | Out-Null
$d.Content.Find.Execute("359×2=", $true, $false, $false, $false, $false, $true, 1, $false, "479×8=", 2) | Out-Null
$d.Content.Find.Execute("770×3=", $true, $false, $false, $false, $false, $true, 1, $false, "154×9=", 2) | Out-Null
$d.Content.Find.Execute("531×9=", $true, $false, $false, $false, $false, $true, 1, $false, "941×9=", 2) | Out-Null
$d.Content.Find.Execute("347×4=", $true, $false, $false, $false, $false, $true, 1, $false, "520×6=", 2) | Out-Null
$d.Content.Find.Execute("264×3=", $true, $false, $false, $false, $false, $true, 1, $false, "449×7=", 2) | Out-Null
$d.Content.Find.Execute("311×5=", $true, $false, $false, $false, $false, $true, 1, $false, "919×6=", 2) | Out-Null
$d.Content.Find.Execute("558×7=", $true, $false, $false, $false, $false, $true, 1, $false, "947×4=", 2) | Out-Null
$d.Content.Find.Execute("589×7=", $true, $false, $false, $false, $false, $true, 1, $false, "450×6=", 2) | Out-Null
$d.Content.Find.Execute("951×3=", $true, $false, $false, $false, $false, $true, 1, $false, "994×6=", 2) | Out-Null

Write-Host "Replacements complete"
